# Replace the placeholder "-" / " " shared-string cells with numeric 0 across
# the various data rows, matching the source workbook's upload.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$cells = @(
    "K8","L8","M8","N8","O8","P8",
    "K9","L9","M9","N9","O9",
    "L11","M11",
    "K14","L14","M14",
    "E15",
    "K15","L15","M15",
    "K17","L17","M17","N17","O17","P17",
    "D25","E25","F25","G25","H25",
    "K25","L25","M25","N25","O25","P25",
    "K33"
)

foreach ($addr in $cells) {
    $ws.Range($addr).Value = 0
}

# Restore the active selection to E16, matching the saved sheet view state.
$ws.Range("E16").Select()
